$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 9
$ws.Cells.Item(9, 1).Value = 112587034
$ws.Cells.Item(9, 2).Value = 89006
$ws.Cells.Item(9, 3).Value = 'Ovaliderad'
$ws.Cells.Item(9, 4).Value = 'LC'
$ws.Cells.Item(9, 5).Value = 4188
$ws.Cells.Item(9, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(9, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(9, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(9, 16).Value = 'Fröjel, Gtl'
$ws.Cells.Item(9, 17).Value = 692906
$ws.Cells.Item(9, 18).Value = 6359283
$ws.Cells.Item(9, 19).Value = 10
$ws.Cells.Item(9, 20).Value = 'Gotland'
$ws.Cells.Item(9, 21).Value = 'Gotland'
$ws.Cells.Item(9, 22).Value = 'Gotland'
$ws.Cells.Item(9, 23).Value = 'Fröjel'
$ws.Cells.Item(9, 25).Value = '''2023-10-07'
$ws.Cells.Item(9, 27).Value = '''2023-10-07'
$ws.Cells.Item(9, 30).Value = $false
$ws.Cells.Item(9, 31).Value = $false
$ws.Cells.Item(9, 33).Value = $false
$ws.Cells.Item(9, 49).Value = 'Brian Johnson'
$ws.Cells.Item(9, 50).Value = 'Brian Johnson, Michael Krikorev, Gillis Aronsson, Helena Björnström'
$ws.Cells.Item(9, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 10
$ws.Cells.Item(10, 1).Value = 112587018
$ws.Cells.Item(10, 2).Value = 89331
$ws.Cells.Item(10, 3).Value = 'Ovaliderad'
$ws.Cells.Item(10, 4).Value = 'LC'
$ws.Cells.Item(10, 5).Value = 3215
$ws.Cells.Item(10, 6).Value = 'Rödgul trumpetsvamp'
$ws.Cells.Item(10, 7).Value = 'Craterellus lutescens'
$ws.Cells.Item(10, 8).Value = '(Fr.) Fr.'
$ws.Cells.Item(10, 16).Value = 'Fröjel, Gtl'
$ws.Cells.Item(10, 17).Value = 692977
$ws.Cells.Item(10, 18).Value = 6359184
$ws.Cells.Item(10, 19).Value = 10
$ws.Cells.Item(10, 20).Value = 'Gotland'
$ws.Cells.Item(10, 21).Value = 'Gotland'
$ws.Cells.Item(10, 22).Value = 'Gotland'
$ws.Cells.Item(10, 23).Value = 'Fröjel'
$ws.Cells.Item(10, 25).Value = '''2023-10-07'
$ws.Cells.Item(10, 27).Value = '''2023-10-07'
$ws.Cells.Item(10, 30).Value = $false
$ws.Cells.Item(10, 31).Value = $false
$ws.Cells.Item(10, 33).Value = $false
$ws.Cells.Item(10, 49).Value = 'Brian Johnson'
$ws.Cells.Item(10, 50).Value = 'Brian Johnson, Michael Krikorev, Gillis Aronsson, Helena Björnström'
$ws.Cells.Item(10, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 11
$ws.Cells.Item(11, 1).Value = 112585364
$ws.Cells.Item(11, 2).Value = 89006
$ws.Cells.Item(11, 3).Value = 'Ovaliderad'
$ws.Cells.Item(11, 4).Value = 'LC'
$ws.Cells.Item(11, 5).Value = 4188
$ws.Cells.Item(11, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(11, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(11, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(11, 16).Value = 'Fröjel-Sigdarve, Gtl'
$ws.Cells.Item(11, 17).Value = 692946
$ws.Cells.Item(11, 18).Value = 6359174
$ws.Cells.Item(11, 19).Value = 10
$ws.Cells.Item(11, 20).Value = 'Gotland'
$ws.Cells.Item(11, 21).Value = 'Gotland'
$ws.Cells.Item(11, 22).Value = 'Gotland'
$ws.Cells.Item(11, 23).Value = 'Fröjel'
$ws.Cells.Item(11, 25).Value = '''2023-10-07'
$ws.Cells.Item(11, 27).Value = '''2023-10-07'
$ws.Cells.Item(11, 30).Value = $false
$ws.Cells.Item(11, 31).Value = $false
$ws.Cells.Item(11, 33).Value = $false
$ws.Cells.Item(11, 49).Value = 'Helena Björnström'
$ws.Cells.Item(11, 50).Value = 'Helena Björnström, Brian Johnson, Michael Krikorev'
$ws.Cells.Item(11, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 12
$ws.Cells.Item(12, 1).Value = 112587033
$ws.Cells.Item(12, 2).Value = 89006
$ws.Cells.Item(12, 3).Value = 'Ovaliderad'
$ws.Cells.Item(12, 4).Value = 'LC'
$ws.Cells.Item(12, 5).Value = 4188
$ws.Cells.Item(12, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(12, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(12, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(12, 16).Value = 'Fröjel, Gtl'
$ws.Cells.Item(12, 17).Value = 692945
$ws.Cells.Item(12, 18).Value = 6359178
$ws.Cells.Item(12, 19).Value = 10
$ws.Cells.Item(12, 20).Value = 'Gotland'
$ws.Cells.Item(12, 21).Value = 'Gotland'
$ws.Cells.Item(12, 22).Value = 'Gotland'
$ws.Cells.Item(12, 23).Value = 'Fröjel'
$ws.Cells.Item(12, 25).Value = '''2023-10-07'
$ws.Cells.Item(12, 27).Value = '''2023-10-07'
$ws.Cells.Item(12, 30).Value = $false
$ws.Cells.Item(12, 31).Value = $false
$ws.Cells.Item(12, 33).Value = $false
$ws.Cells.Item(12, 49).Value = 'Brian Johnson'
$ws.Cells.Item(12, 50).Value = 'Brian Johnson, Michael Krikorev, Gillis Aronsson, Helena Björnström'
$ws.Cells.Item(12, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 13
$ws.Cells.Item(13, 1).Value = 112585365
$ws.Cells.Item(13, 2).Value = 89006
$ws.Cells.Item(13, 3).Value = 'Ovaliderad'
$ws.Cells.Item(13, 4).Value = 'LC'
$ws.Cells.Item(13, 5).Value = 4188
$ws.Cells.Item(13, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(13, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(13, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(13, 16).Value = 'Fröjel-Sigdarve, Gtl'
$ws.Cells.Item(13, 17).Value = 692959
$ws.Cells.Item(13, 18).Value = 6359217
$ws.Cells.Item(13, 19).Value = 10
$ws.Cells.Item(13, 20).Value = 'Gotland'
$ws.Cells.Item(13, 21).Value = 'Gotland'
$ws.Cells.Item(13, 22).Value = 'Gotland'
$ws.Cells.Item(13, 23).Value = 'Fröjel'
$ws.Cells.Item(13, 25).Value = '''2023-10-07'
$ws.Cells.Item(13, 27).Value = '''2023-10-07'
$ws.Cells.Item(13, 30).Value = $false
$ws.Cells.Item(13, 31).Value = $false
$ws.Cells.Item(13, 33).Value = $false
$ws.Cells.Item(13, 49).Value = 'Helena Björnström'
$ws.Cells.Item(13, 50).Value = 'Helena Björnström, Brian Johnson, Michael Krikorev'
$ws.Cells.Item(13, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 14
$ws.Cells.Item(14, 1).Value = 112586310
$ws.Cells.Item(14, 2).Value = 89336
$ws.Cells.Item(14, 3).Value = 'Ovaliderad'
$ws.Cells.Item(14, 4).Value = 'VU'
$ws.Cells.Item(14, 5).Value = 2015
$ws.Cells.Item(14, 6).Value = 'Vit taggsvamp'
$ws.Cells.Item(14, 7).Value = 'Hydnum albidum'
$ws.Cells.Item(14, 8).Value = 'Peck'
$ws.Cells.Item(14, 9).Value = '''2'
$ws.Cells.Item(14, 10).Value = 'fruktkroppar'
$ws.Cells.Item(14, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(14, 17).Value = 692897
$ws.Cells.Item(14, 18).Value = 6359273
$ws.Cells.Item(14, 19).Value = 10
$ws.Cells.Item(14, 20).Value = 'Gotland'
$ws.Cells.Item(14, 21).Value = 'Gotland'
$ws.Cells.Item(14, 22).Value = 'Gotland'
$ws.Cells.Item(14, 23).Value = 'Fröjel'
$ws.Cells.Item(14, 25).Value = '''2023-10-07'
$ws.Cells.Item(14, 27).Value = '''2023-10-07'
$ws.Cells.Item(14, 29).Value = '2 ex.'
$ws.Cells.Item(14, 30).Value = $false
$ws.Cells.Item(14, 31).Value = $false
$ws.Cells.Item(14, 33).Value = $false
$ws.Cells.Item(14, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(14, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(14, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(14, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 15
$ws.Cells.Item(15, 1).Value = 112586312
$ws.Cells.Item(15, 2).Value = 89336
$ws.Cells.Item(15, 3).Value = 'Ovaliderad'
$ws.Cells.Item(15, 4).Value = 'VU'
$ws.Cells.Item(15, 5).Value = 2015
$ws.Cells.Item(15, 6).Value = 'Vit taggsvamp'
$ws.Cells.Item(15, 7).Value = 'Hydnum albidum'
$ws.Cells.Item(15, 8).Value = 'Peck'
$ws.Cells.Item(15, 9).Value = '''5'
$ws.Cells.Item(15, 10).Value = 'fruktkroppar'
$ws.Cells.Item(15, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(15, 17).Value = 692883
$ws.Cells.Item(15, 18).Value = 6359288
$ws.Cells.Item(15, 19).Value = 10
$ws.Cells.Item(15, 20).Value = 'Gotland'
$ws.Cells.Item(15, 21).Value = 'Gotland'
$ws.Cells.Item(15, 22).Value = 'Gotland'
$ws.Cells.Item(15, 23).Value = 'Fröjel'
$ws.Cells.Item(15, 25).Value = '''2023-10-07'
$ws.Cells.Item(15, 27).Value = '''2023-10-07'
$ws.Cells.Item(15, 29).Value = 'Ca. 5 ex.'
$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false
$ws.Cells.Item(15, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(15, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(15, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(15, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 16
$ws.Cells.Item(16, 1).Value = 112586311
$ws.Cells.Item(16, 2).Value = 89006
$ws.Cells.Item(16, 3).Value = 'Ovaliderad'
$ws.Cells.Item(16, 4).Value = 'LC'
$ws.Cells.Item(16, 5).Value = 4188
$ws.Cells.Item(16, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(16, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(16, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(16, 9).Value = '''2'
$ws.Cells.Item(16, 10).Value = 'fruktkroppar'
$ws.Cells.Item(16, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(16, 17).Value = 692886
$ws.Cells.Item(16, 18).Value = 6359280
$ws.Cells.Item(16, 19).Value = 10
$ws.Cells.Item(16, 20).Value = 'Gotland'
$ws.Cells.Item(16, 21).Value = 'Gotland'
$ws.Cells.Item(16, 22).Value = 'Gotland'
$ws.Cells.Item(16, 23).Value = 'Fröjel'
$ws.Cells.Item(16, 25).Value = '''2023-10-07'
$ws.Cells.Item(16, 27).Value = '''2023-10-07'
$ws.Cells.Item(16, 29).Value = '2 ex. i barrmatta under gran.'
$ws.Cells.Item(16, 30).Value = $false
$ws.Cells.Item(16, 31).Value = $false
$ws.Cells.Item(16, 33).Value = $false
$ws.Cells.Item(16, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(16, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(16, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(16, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 17
$ws.Cells.Item(17, 1).Value = 112586308
$ws.Cells.Item(17, 2).Value = 89006
$ws.Cells.Item(17, 3).Value = 'Ovaliderad'
$ws.Cells.Item(17, 4).Value = 'LC'
$ws.Cells.Item(17, 5).Value = 4188
$ws.Cells.Item(17, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(17, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(17, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(17, 9).Value = '''6'
$ws.Cells.Item(17, 10).Value = 'fruktkroppar'
$ws.Cells.Item(17, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(17, 17).Value = 692968
$ws.Cells.Item(17, 18).Value = 6359239
$ws.Cells.Item(17, 19).Value = 10
$ws.Cells.Item(17, 20).Value = 'Gotland'
$ws.Cells.Item(17, 21).Value = 'Gotland'
$ws.Cells.Item(17, 22).Value = 'Gotland'
$ws.Cells.Item(17, 23).Value = 'Fröjel'
$ws.Cells.Item(17, 25).Value = '''2023-10-07'
$ws.Cells.Item(17, 27).Value = '''2023-10-07'
$ws.Cells.Item(17, 29).Value = '6 ex. i barrmatta under gran.'
$ws.Cells.Item(17, 30).Value = $false
$ws.Cells.Item(17, 31).Value = $false
$ws.Cells.Item(17, 33).Value = $false
$ws.Cells.Item(17, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(17, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(17, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(17, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 18
$ws.Cells.Item(18, 1).Value = 112586309
$ws.Cells.Item(18, 2).Value = 89006
$ws.Cells.Item(18, 3).Value = 'Ovaliderad'
$ws.Cells.Item(18, 4).Value = 'LC'
$ws.Cells.Item(18, 5).Value = 4188
$ws.Cells.Item(18, 6).Value = 'Fransig jordstjärna'
$ws.Cells.Item(18, 7).Value = 'Geastrum fimbriatum'
$ws.Cells.Item(18, 8).Value = 'Fr.:Pers.'
$ws.Cells.Item(18, 9).Value = '''1'
$ws.Cells.Item(18, 10).Value = 'fruktkroppar'
$ws.Cells.Item(18, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(18, 17).Value = 692961
$ws.Cells.Item(18, 18).Value = 6359243
$ws.Cells.Item(18, 19).Value = 10
$ws.Cells.Item(18, 20).Value = 'Gotland'
$ws.Cells.Item(18, 21).Value = 'Gotland'
$ws.Cells.Item(18, 22).Value = 'Gotland'
$ws.Cells.Item(18, 23).Value = 'Fröjel'
$ws.Cells.Item(18, 25).Value = '''2023-10-07'
$ws.Cells.Item(18, 27).Value = '''2023-10-07'
$ws.Cells.Item(18, 29).Value = '1 ex. i barrmatta under gran.'
$ws.Cells.Item(18, 30).Value = $false
$ws.Cells.Item(18, 31).Value = $false
$ws.Cells.Item(18, 33).Value = $false
$ws.Cells.Item(18, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(18, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(18, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(18, 51).Value = 'Kalkbarrianerna 2023, Gotland'

# Row 19
$ws.Cells.Item(19, 1).Value = 112586306
$ws.Cells.Item(19, 2).Value = 89336
$ws.Cells.Item(19, 3).Value = 'Ovaliderad'
$ws.Cells.Item(19, 4).Value = 'VU'
$ws.Cells.Item(19, 5).Value = 2015
$ws.Cells.Item(19, 6).Value = 'Vit taggsvamp'
$ws.Cells.Item(19, 7).Value = 'Hydnum albidum'
$ws.Cells.Item(19, 8).Value = 'Peck'
$ws.Cells.Item(19, 9).Value = '''1'
$ws.Cells.Item(19, 10).Value = 'fruktkroppar'
$ws.Cells.Item(19, 16).Value = 'Syrmansberget vid Fröjel, Gtl'
$ws.Cells.Item(19, 17).Value = 692937
$ws.Cells.Item(19, 18).Value = 6359174
$ws.Cells.Item(19, 19).Value = 10
$ws.Cells.Item(19, 20).Value = 'Gotland'
$ws.Cells.Item(19, 21).Value = 'Gotland'
$ws.Cells.Item(19, 22).Value = 'Gotland'
$ws.Cells.Item(19, 23).Value = 'Fröjel'
$ws.Cells.Item(19, 25).Value = '''2023-10-07'
$ws.Cells.Item(19, 27).Value = '''2023-10-07'
$ws.Cells.Item(19, 29).Value = '1 ex.'
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 35).Value = 'Äldre kalktallskog med inslag av gran.'
$ws.Cells.Item(19, 49).Value = 'Gillis Aronsson'
$ws.Cells.Item(19, 50).Value = 'Gillis Aronsson, Helena Björnström, Brian Johnson'
$ws.Cells.Item(19, 51).Value = 'Kalkbarrianerna 2023, Gotland'

